$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1666666666666667
$ws.Range("C2").Value = 0.5944444444444444
$ws.Range("J2").Value = 0.002777777777777778
$ws.Range("P2").Value = 0.1527777777777778
$ws.Range("S2").Value = 0.08333333333333333
$ws.Range("C3").Value = 0.02764976958525346
$ws.Range("J3").Value = 0.05069124423963134
$ws.Range("P3").Value = 0.7465437788018433
$ws.Range("S3").Value = 0.1751152073732719
$ws.Range("J4").Value = 0.03448275862068965
$ws.Range("P4").Value = 0.7758620689655172
$ws.Range("S4").Value = 0.1896551724137931
$ws.Range("B6").Value = 0.09433962264150944
$ws.Range("D6").Value = 0.009433962264150943
$ws.Range("F6").Value = 0.04716981132075472
$ws.Range("J6").Value = 0.2216981132075472
$ws.Range("O6").Value = 0.0330188679245283
$ws.Range("Q6").Value = 0.2122641509433962
$ws.Range("R6").Value = 0.02358490566037736
$ws.Range("S6").Value = 0.3584905660377358
$ws.Range("B7").Value = 0.1641025641025641
$ws.Range("D7").Value = 0.01538461538461539
$ws.Range("E7").Value = 0.005128205128205128
$ws.Range("F7").Value = 0.03589743589743589
$ws.Range("J7").Value = 0.1230769230769231
$ws.Range("O7").Value = 0.01538461538461539
$ws.Range("Q7").Value = 0.1230769230769231
$ws.Range("R7").Value = 0.07179487179487179
$ws.Range("S7").Value = 0.4461538461538462
$ws.Range("B8").Value = 0.1118721461187215
$ws.Range("D8").Value = 0.0273972602739726
$ws.Range("F8").Value = 0.0639269406392694
$ws.Range("J8").Value = 0.07990867579908675
$ws.Range("O8").Value = 0.0136986301369863
$ws.Range("Q8").Value = 0.1963470319634703
$ws.Range("R8").Value = 0.08447488584474885
$ws.Range("S8").Value = 0.4223744292237443
$ws.Range("B9").Value = 0.1078838174273859
$ws.Range("D9").Value = 0.02074688796680498
$ws.Range("E9").Value = 0.004149377593360996
$ws.Range("F9").Value = 0.06639004149377593
$ws.Range("J9").Value = 0.06639004149377593
$ws.Range("O9").Value = 0.01659751037344398
$ws.Range("Q9").Value = 0.1991701244813278
$ws.Range("R9").Value = 0.07468879668049792
$ws.Range("S9").Value = 0.4439834024896265
$ws.Range("B10").Value = 0.1290064102564103
$ws.Range("D10").Value = 0.03044871794871795
$ws.Range("E10").Value = 0.001602564102564103
$ws.Range("F10").Value = 0.06490384615384616
$ws.Range("J10").Value = 0.1017628205128205
$ws.Range("O10").Value = 0.01923076923076923
$ws.Range("Q10").Value = 0.2355769230769231
$ws.Range("R10").Value = 0.08253205128205128
$ws.Range("S10").Value = 0.3349358974358974
$ws.Range("G11").Value = 0.1413793103448276
$ws.Range("J11").Value = 0.06551724137931035
$ws.Range("K11").Value = 0.1724137931034483
$ws.Range("L11").Value = 0.5862068965517241
$ws.Range("S11").Value = 0.03448275862068965
$ws.Range("G12").Value = 0.7336956521739131
$ws.Range("J12").Value = 0.1793478260869565
$ws.Range("L12").Value = 0.04347826086956522
$ws.Range("S12").Value = 0.04347826086956522
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.3055555555555556
$ws.Range("S13").Value = 0.02777777777777778
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.01646090534979424
$ws.Range("H15").Value = 0.1769547325102881
$ws.Range("I15").Value = 0.07818930041152264
$ws.Range("J15").Value = 0.3580246913580247
$ws.Range("K15").Value = 0.04938271604938271
$ws.Range("M15").Value = 0.00411522633744856
$ws.Range("O15").Value = 0.06584362139917696
$ws.Range("S15").Value = 0.2510288065843622
$ws.Range("F16").Value = 0.003846153846153846
$ws.Range("H16").Value = 0.1769230769230769
$ws.Range("I16").Value = 0.1
$ws.Range("J16").Value = 0.4423076923076923
$ws.Range("K16").Value = 0.09615384615384616
$ws.Range("M16").Value = 0.02307692307692308
$ws.Range("N16").Value = 0.003846153846153846
$ws.Range("O16").Value = 0.03461538461538462
$ws.Range("S16").Value = 0.1192307692307692
$ws.Range("F17").Value = 0.01609657947686117
$ws.Range("H17").Value = 0.1670020120724346
$ws.Range("I17").Value = 0.1207243460764588
$ws.Range("J17").Value = 0.4044265593561368
$ws.Range("K17").Value = 0.1146881287726358
$ws.Range("M17").Value = 0.01006036217303823
$ws.Range("O17").Value = 0.07847082494969819
$ws.Range("S17").Value = 0.08853118712273642
$ws.Range("F18").Value = 0.02906976744186046
$ws.Range("H18").Value = 0.1918604651162791
$ws.Range("I18").Value = 0.1046511627906977
$ws.Range("J18").Value = 0.4069767441860465
$ws.Range("K18").Value = 0.1162790697674419
$ws.Range("M18").Value = 0.02325581395348837
$ws.Range("O18").Value = 0.05813953488372093
$ws.Range("S18").Value = 0.06976744186046512
$ws.Range("F19").Value = 0.02333065164923572
$ws.Range("H19").Value = 0.1922767497988737
$ws.Range("I19").Value = 0.09814963797264682
$ws.Range("J19").Value = 0.3716814159292036
$ws.Range("K19").Value = 0.09734513274336283
$ws.Range("M19").Value = 0.01609010458567981
$ws.Range("N19").Value = 0.0008045052292839903
$ws.Range("O19").Value = 0.07401448109412712
$ws.Range("S19").Value = 0.1263073209975865
